# Removed Test Case Inter-Dependency
# Updates the product name to an independent, suffixed value and the
# shortname to a non-numeric (string) identifier so this test case no
# longer collides with / depends on other test cases that reuse "4353".

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# productname on both the input and output sheets
$wsInput.Range("B1").Value = "4353-MS-Simple-Group-Loan-Product-Loanproduct-OVERDUEFEEFLAT-1st"
$wsOutput.Range("B1").Value = "4353-MS-Simple-Group-Loan-Product-Loanproduct-OVERDUEFEEFLAT-1st"

# shortname becomes a text value instead of a plain number
$wsInput.Range("B2").Value = "435p"

# Move the active selection (cosmetic) to B3
$wsInput.Activate()
$wsInput.Range("B3").Select()
